$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-117 down to 116-118
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly record
$ws.Range("A115").Value = 6
$ws.Range("B115").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C115").Value = "Metropolitana"
$ws.Range("D115").Value = 45008
$ws.Range("E115").Value = 13
$ws.Range("F115").Value = 100114007
$ws.Range("G115").Value = "Jengibre"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 400
$ws.Range("K115").Value = 15000
$ws.Range("L115").Value = 16000
$ws.Range("M115").Value = 15425
$ws.Range("N115").Value = "$/caja 13 kilos"
$ws.Range("O115").Value = "Perú"
$ws.Range("P115").Value = 1187
$ws.Range("Q115").Value = 13
$ws.Range("R115").Value = "Hortaliza"
